$wb = $excel.ActiveWorkbook

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1012095.75
$ws.Range("I17").Value = 622
$ws.Range("J17").Value = 1113243.1
$ws.Range("K17").Value = 1866
$ws.Range("L17").Value = 3339729.3
$ws.Range("M17").Value = -1698
$ws.Range("N17").Value = -3340065.3

# ALC row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3724.9167
$ws.Range("I43").Value = 4583.1665
$ws.Range("J43").Value = 2866.6667
$ws.Range("K43").Value = 4583.1665
$ws.Range("L43").Value = 2866.6667
$ws.Range("M43").Value = -4514.1665
$ws.Range("N43").Value = -3004.6667

# ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4450.8125
$ws.Range("I74").Value = 2459
$ws.Range("K74").Value = 2459
$ws.Range("M74").Value = -1523

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 13173.714
$ws.Range("I76").Value = 29120.5
$ws.Range("J76").Value = 6795
$ws.Range("K76").Value = 29120.5
$ws.Range("L76").Value = 6795
$ws.Range("M76").Value = -28805.5
$ws.Range("N76").Value = -7425

# ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 4450.8125
$ws.Range("I77").Value = 2459
$ws.Range("K77").Value = 12295
$ws.Range("M77").Value = -7615

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 13173.714
$ws.Range("I79").Value = 29120.5
$ws.Range("J79").Value = 6795
$ws.Range("K79").Value = 29120.5
$ws.Range("L79").Value = 6795
$ws.Range("M79").Value = -28028.5
$ws.Range("N79").Value = -8979

# ALC row 101
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 1236
$ws.Range("I101").Value = 1314.8334
$ws.Range("J101").Value = 999.5
$ws.Range("K101").Value = 3944.5002
$ws.Range("L101").Value = 2998.5
$ws.Range("M101").Value = -2322.5002
$ws.Range("N101").Value = -6242.5

# ALC row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 556.28
$ws.Range("I107").Value = 534.25
$ws.Range("J107").Value = 644.4
$ws.Range("K107").Value = 534.25
$ws.Range("L107").Value = 644.4
$ws.Range("M107").Value = 1385.75
$ws.Range("N107").Value = -4484.4

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 4095.125
$ws.Range("I113").Value = 1830.8
$ws.Range("J113").Value = 5124.364
$ws.Range("K113").Value = 1830.8
$ws.Range("L113").Value = 5124.364
$ws.Range("M113").Value = 1423.2
$ws.Range("N113").Value = -11632.364

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4866.25
$ws.Range("I132").Value = 5159.375
$ws.Range("K132").Value = 15478.125
$ws.Range("M132").Value = -12948.125

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16215.542
$ws.Range("I32").Value = 13840.667
$ws.Range("K32").Value = 13840.667
$ws.Range("M32").Value = -13553.667

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2121.5
$ws.Range("I45").Value = 1824
$ws.Range("K45").Value = 1824
$ws.Range("M45").Value = -1447

# ARM row 46
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 8664.333000000001
$ws.Range("I46").Value = 8995
$ws.Range("J46").Value = 8499
$ws.Range("K46").Value = 8995
$ws.Range("L46").Value = 8499
$ws.Range("M46").Value = -8676
$ws.Range("N46").Value = -9137

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 15629986
$ws.Range("I61").Value = 22730758
$ws.Range("K61").Value = 22730758
$ws.Range("M61").Value = -22730546

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 8500.625
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 8500.625
$ws.Range("K63").Value = 0
$ws.Range("L63").ClearContents()
$ws.Range("M63").Value = 8500.625
$ws.Range("N63").Value = -9872.625

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 8500.625
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 8500.625
$ws.Range("K66").Value = 0
$ws.Range("L66").ClearContents()
$ws.Range("M66").Value = 42503.125
$ws.Range("N66").Value = -49367.125

# ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 15750.714
$ws.Range("I88").Value = 50553
$ws.Range("J88").Value = 1829.8
$ws.Range("K88").Value = 50553
$ws.Range("L88").Value = 1829.8
$ws.Range("M88").Value = -50147
$ws.Range("N88").Value = -2641.8

# ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 15750.714
$ws.Range("I91").Value = 50553
$ws.Range("J91").Value = 1829.8
$ws.Range("K91").Value = 50553
$ws.Range("L91").Value = 1829.8
$ws.Range("M91").Value = -49149
$ws.Range("N91").Value = -4637.8

# ARM row 113
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H113").Value = 45000
$ws.Range("J113").Value = 45000
$ws.Range("L113").Value = 45000
$ws.Range("N113").Value = -53678

# ARM row 114
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H114").Value = 50599.5
$ws.Range("J114").Value = 50599.5
$ws.Range("L114").Value = 50599.5
$ws.Range("N114").Value = -59277.5

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 35774332
$ws.Range("I132").Value = 11440.625
$ws.Range("K132").Value = 34321.875
$ws.Range("M132").Value = -31791.875

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 15629986
$ws.Range("I136").Value = 22730758
$ws.Range("K136").Value = 68192274
$ws.Range("M136").Value = -68189724

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2189
$ws.Range("J20").Value = 2000
$ws.Range("L20").Value = 2000
$ws.Range("N20").Value = -2494

# BSM row 80
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 968.05
$ws.Range("I80").Value = 597.5
$ws.Range("J80").Value = 1523.875
$ws.Range("K80").Value = 597.5
$ws.Range("L80").Value = 1523.875
$ws.Range("M80").Value = 400.5
$ws.Range("N80").Value = -3519.875

# BSM row 83
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 968.05
$ws.Range("I83").Value = 597.5
$ws.Range("J83").Value = 1523.875
$ws.Range("K83").Value = 2987.5
$ws.Range("L83").Value = 7619.375
$ws.Range("M83").Value = 2004.5
$ws.Range("N83").Value = -17603.375

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 20071.389
$ws.Range("I86").Value = 10328.826
$ws.Range("K86").Value = 10328.826
$ws.Range("M86").Value = -9205.825999999999

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 20071.389
$ws.Range("I89").Value = 10328.826
$ws.Range("K89").Value = 51644.13
$ws.Range("M89").Value = -46028.13

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2739.4348
$ws.Range("I134").Value = 2700.5454
$ws.Range("K134").Value = 8101.6362
$ws.Range("M134").Value = -5566.6362

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 60059.484
$ws.Range("I132").Value = 85491.586
$ws.Range("K132").Value = 256474.758
$ws.Range("M132").Value = -253944.758

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2279.3809
$ws.Range("I134").Value = 2051.0588
$ws.Range("K134").Value = 6153.176399999999
$ws.Range("M134").Value = -3618.176399999999

# CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 28153842
$ws.Range("I4").Value = 35187010
$ws.Range("J4").Value = 17955750
$ws.Range("K4").Value = 105561030
$ws.Range("L4").Value = 53867250
$ws.Range("M4").Value = -105560918
$ws.Range("N4").Value = -53867474

# CUL row 130
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 2531.5
$ws.Range("J130").Value = 3033
$ws.Range("L130").Value = 9099
$ws.Range("N130").Value = -19139

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1667.5483
$ws.Range("J131").Value = 1719.4286
$ws.Range("L131").Value = 5158.2858
$ws.Range("N131").Value = -15238.2858

# CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1010.34784
$ws.Range("J140").Value = 3032
$ws.Range("L140").Value = 9096
$ws.Range("N140").Value = -19456

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 1002504
$ws.Range("I70").Value = 2000008
$ws.Range("K70").Value = 2000008
$ws.Range("M70").Value = -1999738

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 1002504
$ws.Range("I73").Value = 2000008
$ws.Range("K73").Value = 2000008
$ws.Range("M73").Value = -1999072

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3153.25
$ws.Range("I80").Value = 2966.4285
$ws.Range("K80").Value = 2966.4285
$ws.Range("M80").Value = -1968.4285

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3153.25
$ws.Range("I83").Value = 2966.4285
$ws.Range("K83").Value = 14832.1425
$ws.Range("M83").Value = -9840.1425

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2204
$ws.Range("I132").Value = 1806
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 5418
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -2888
$ws.Range("N132").Value = -14060

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4963.3335
$ws.Range("J68").Value = 6050
$ws.Range("L68").Value = 6050
$ws.Range("N68").Value = -7548

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 4963.3335
$ws.Range("J71").Value = 6050
$ws.Range("L71").Value = 30250
$ws.Range("N71").Value = -37738

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5041.9116
$ws.Range("I122").Value = 3988
$ws.Range("K122").Value = 11964
$ws.Range("M122").Value = -9514

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 32470
$ws.Range("I132").Value = 65997.5
$ws.Range("J132").Value = 22153.846
$ws.Range("K132").Value = 197992.5
$ws.Range("L132").Value = 66461.538
$ws.Range("M132").Value = -195462.5
$ws.Range("N132").Value = -71521.538

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2200.139
$ws.Range("I136").Value = 1022.64
$ws.Range("J136").Value = 4876.273
$ws.Range("K136").Value = 3067.92
$ws.Range("L136").Value = 14628.819
$ws.Range("M136").Value = -517.9200000000001
$ws.Range("N136").Value = -19728.819

# WVR row 14
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 778.2222
$ws.Range("I14").Value = 500
$ws.Range("J14").Value = 1334.6666
$ws.Range("K14").Value = 500
$ws.Range("L14").Value = 1334.6666
$ws.Range("M14").Value = -332
$ws.Range("N14").Value = -1670.6666

# WVR row 31
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 14998
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7428.5713
$ws.Range("I62").Value = 6000
$ws.Range("K62").Value = 6000
$ws.Range("M62").Value = -5376

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 7428.5713
$ws.Range("I65").Value = 6000
$ws.Range("K65").Value = 30000
$ws.Range("M65").Value = -26880

# WVR row 123
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 69999
$ws.Range("J123").Value = 69999
$ws.Range("L123").Value = 69999
$ws.Range("N123").Value = -79799

# WVR row 133
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 100000
$ws.Range("J133").Value = 100000
$ws.Range("L133").Value = 100000
$ws.Range("N133").Value = -110120
